# Updates cryptocurrency price (column D) and 1h volume-change (column E)
# figures on the active worksheet, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'63.514.42"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Formula = "'2.637.55"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Formula = "'605.71"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Formula = "'147.43"
$ws.Range("E6").Value = "  +3.29%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Formula = "'0.590"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").Formula = "'0.373"
$ws.Range("E11").Value = "  +4.96%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Formula = "'27.64"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Formula = "'3.108.94"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Formula = "'63.345.54"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Formula = "'2.636.67"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Formula = "'11.58"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").Formula = "'4.59"
$ws.Range("E19").Value = "  +5.44%  "
$ws.Range("D20").Formula = "'344.53"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("D22").Formula = "'0.999"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D24").Formula = "'66.88"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("D26").Formula = "'9.08"
$ws.Range("E26").Value = "  +8.35%  "
$ws.Range("D27").Formula = "'1.56"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Formula = "'566.04"
$ws.Range("E28").Value = "  +7.74%  "
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").Formula = "'8.00"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("E33").Value = "  +6.45%  "
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("D35").Formula = "'5.17"
$ws.Range("E35").Value = "  +6.01%  "
$ws.Range("D36").Formula = "'167.15"
$ws.Range("E36").Value = "  -4.16%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").Formula = "'1.00"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +8.74%  "
$ws.Range("D40").Formula = "'19.18"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Formula = "'166.95"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Formula = "'3.80"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").Formula = "'22.23"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Formula = "'0.0572"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Formula = "'0.0247"
$ws.Range("E47").Value = "  +4.09%  "
$ws.Range("D48").Formula = "'0.0963"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Formula = "'1.91"
$ws.Range("E49").Value = "  +13.47%  "
$ws.Range("D50").Formula = "'18.86"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Formula = "'0.183"
$ws.Range("E51").Value = "  +5.91%  "
